$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-17 11:02:27"
$wsDeDe.Range("H2").Value = "2016-08-17 11:02:27"
$wsZhCn.Range("H2").Value = "2016-08-17 11:02:22"
$wsZhCn.Range("K2").Value = "2016-08-17 11:02:40"
$wsDeDe.Range("K2").Value = "2016-08-17 11:02:47"
